# Insert two new rows at row 482, pushing the existing rows (482:546) down to (484:548).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("482:483").Insert()

# New row 482: Calidad "Primera", fecha 44918 (2022-12-23)
$ws.Range("A482").Value = 8
$ws.Range("B482").Value = "Terminal La Palmera de La Serena"
$ws.Range("C482").Value = "Coquimbo"
$ws.Range("D482").Value = 44918
$ws.Range("E482").Value = 4
$ws.Range("F482").Value = 100112009
$ws.Range("G482").Value = "Acelga"
$ws.Range("H482").Value = "Sin especificar"
$ws.Range("I482").Value = "Primera"
$ws.Range("J482").Value = 2000
$ws.Range("K482").Value = 600
$ws.Range("L482").Value = 700
$ws.Range("M482").Value = 650
$ws.Range("N482").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O482").Value = "Provincia del Elqu$([char]0xED)"
$ws.Range("P482").Value = 325
$ws.Range("Q482").Value = 2
$ws.Range("R482").Value = "Hortaliza"

# New row 483: Calidad "Segunda", fecha 44918 (2022-12-23)
$ws.Range("A483").Value = 8
$ws.Range("B483").Value = "Terminal La Palmera de La Serena"
$ws.Range("C483").Value = "Coquimbo"
$ws.Range("D483").Value = 44918
$ws.Range("E483").Value = 4
$ws.Range("F483").Value = 100112009
$ws.Range("G483").Value = "Acelga"
$ws.Range("H483").Value = "Sin especificar"
$ws.Range("I483").Value = "Segunda"
$ws.Range("J483").Value = 1520
$ws.Range("K483").Value = 450
$ws.Range("L483").Value = 500
$ws.Range("M483").Value = 475
$ws.Range("N483").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O483").Value = "Provincia del Elqu$([char]0xED)"
$ws.Range("P483").Value = 238
$ws.Range("Q483").Value = 2
$ws.Range("R483").Value = "Hortaliza"
